# Apply crypto price/volume updates per the commit diff (Mon Nov 4 19:42:59 UTC 2024).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.923.41'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.01%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.433.79'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.63%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '160.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.66%  '

# Row 7
$ws.Range("E7").Value = '  +0.10%  '

# Row 8
$ws.Range("E8").Value = '  +1.75%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.160'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +7.40%  '

# Row 10
$ws.Range("E10").Value = '  -0.61%  '

# Row 11
$ws.Range("E11").Value = '  -1.25%  '

# Row 12
$ws.Range("E12").Value = '  -0.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '67.830.78'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.86%  '

# Row 14
$ws.Range("E14").Value = '  +1.47%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.01'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.15%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '10.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.19%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '334.09'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.73%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.90%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.79'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.23%  '

# Row 20
$ws.Range("E20").Value = '  +0.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.43%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.16'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.30%  '

# Row 23
$ws.Range("E23").Value = '  -0.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.11'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.47%  '

# Row 25
$ws.Range("E25").Value = '  -0.39%  '

# Row 26
$ws.Range("E26").Value = '  -0.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.04%  '

# Row 28
$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.76%  '

# Row 29
$ws.Range("B29").Value = 'Bittensor'
$ws.Range("C29").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '417.06'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.93%  '

# Row 30
$ws.Range("E30").Value = '  -0.27%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '160.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.86%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.95'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.53%  '

# Row 33
$ws.Range("E33").Value = '  -0.11%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '17.85'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.21%  '

# Row 35
$ws.Range("E35").Value = '  -3.32%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.295'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.21%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.26'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.02%  '

# Row 38
$ws.Range("E38").Value = '  +0.89%  '

# Row 39
$ws.Range("E39").Value = '  -1.01%  '

# Row 40
$ws.Range("E40").Value = '  -1.74%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.33'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '129.57'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.64%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0706'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.74%  '

# Row 44
$ws.Range("E44").Value = '  -0.38%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.556'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.42%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0914'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.92%  '

# Row 47
$ws.Range("E47").Value = '  +0.55%  '

# Row 48
$ws.Range("E48").Value = '  -6.74%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '16.54'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.40%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0203'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.26%  '

# Row 51
$ws.Range("E51").Value = '  +0.52%  '
